# Bulk upload template fix:
#  - Rename the "ResourceCreator.Institution.name" / "ResourceCreator.role"
#    headers (columns AC/AD) to "ResourceCreatorInstitution.Institution.name"
#    / "ResourceCreatorInstitution.role" so the header text matches the bean
#    property path used by the (newly extracted) ResourceCreatorInstitution
#    bean.
#  - Re-fit the view: widen column G (copyLocation) to fit its header,
#    re-size the header/data rows for the new text and select the new last
#    header cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Header text fix (columns AC/AD on row 1) -----------------------------
$ws.Range("AC1").Value = "ResourceCreatorInstitution.Institution.name"
$ws.Range("AD1").Value = "ResourceCreatorInstitution.role"

# --- Column sizing (bestFit-like width on column G) -----------------------
$ws.Columns.Item(7).ColumnWidth = 32.5

# --- Row heights (re-wrapped header/data rows) -----------------------------
$ws.Rows.Item(1).RowHeight = 51.75
$ws.Rows.Item(2).RowHeight = 75
$ws.Rows.Item(3).RowHeight = 105

# --- View state: scroll right and select the new last header cell --------
$ws.Activate() | Out-Null
$ws.Range("AD1").Select() | Out-Null
try {
    $excel.ActiveWindow.ScrollColumn = 7
} catch {
}
